# Applies the "feat: add 2022-Q1 data" change:
#  - Adds a new sheet "2022-Q1" (positioned right before "总计") holding the
#    per-fund holdings detail for the new quarter.
#  - Updates the "总计" (summary) sheet with a new leading row for 2022-Q1
#    and renumbers the index column.

$wb = $excel.ActiveWorkbook

# Reference cells whose formatting (bold + border + centered header /
# column-A look) we will reuse for the newly written cells, so we don't
# have to hand roll a brand-new style.
$refSheet = $wb.Worksheets.Item("2021-Q4")
$refHeaderCell = $refSheet.Range("B1")
$refIndexCell = $refSheet.Range("A2")

# ------------------------------------------------------------------
# Step 1: create the sheets.
# The existing "总计" worksheet is renamed to "2022-Q1" and repurposed to
# hold the new quarter detail, while a freshly added worksheet becomes the
# new "总计" summary sheet (placed right after the renamed one, i.e. in the
# same slot "总计" used to occupy).
# ------------------------------------------------------------------
$oldTotalSheet = $wb.Worksheets.Item("总计")
$newTotalSheet = $wb.Worksheets.Add($null, $oldTotalSheet)
$newTotalSheet.Name = "总计_tmp"
$oldTotalSheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# ------------------------------------------------------------------
# Step 2: populate the "2022-Q1" detail sheet.
# ------------------------------------------------------------------
$q1Sheet.Cells.Clear()

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

$q1Data = @(
    @("005821", "万家新机遇龙头企业灵活配置混合", "23.23", "56.20", "3.31", "0.7689", 5),
    @("013960", "万家新机遇成长一年持有期混合A", "13.29", "49.51", "3.40", "0.4519", 3),
    @("168501", "北信瑞丰产业升级多策略混合",   "4.42",  "94.11", "3.35", "0.1481", 10),
    @("013961", "万家新机遇成长一年持有期混合C", "3.13",  "49.51", "3.40", "0.1064", 3),
    @("001056", "北信瑞丰健康生活主题灵活配置混合", "1.64", "86.03", "3.61", "0.0592", 8),
    @("519677", "银河定投宝中证腾讯济安价值100A股指数", "2.74", "91.56", "1.27", "0.0348", 7),
    @("005569", "中融智选红利股票A", "0.27", "92.46", "3.87", "0.0104", 6),
    @("001866", "北信瑞丰新成长灵活配置混合", "0.07", "94.21", "4.15", "0.0029", 8),
    @("005570", "中融智选红利股票C", "0.03", "92.46", "3.87", "0.0012", 6)
)

$lastRow = 1 + $q1Data.Count

# Force text storage (keeps leading zeros on fund codes, keeps the numeric
# looking strings such as "23.23" as literal text instead of floats).
$q1Sheet.Range("B2:G" + $lastRow).NumberFormat = "@"

$row = 2
foreach ($item in $q1Data) {
    $q1Sheet.Cells.Item($row, 1).Value = $row - 2
    $q1Sheet.Cells.Item($row, 2).Value = $item[0]
    $q1Sheet.Cells.Item($row, 3).Value = $item[1]
    $q1Sheet.Cells.Item($row, 4).Value = $item[2]
    $q1Sheet.Cells.Item($row, 5).Value = $item[3]
    $q1Sheet.Cells.Item($row, 6).Value = $item[4]
    $q1Sheet.Cells.Item($row, 7).Value = $item[5]
    $q1Sheet.Cells.Item($row, 8).Value = $item[6]
    $row = $row + 1
}

# Match the look of the other quarterly sheets: bold/bordered/centered
# header row and index (column A) cells.
$refHeaderCell.Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$refIndexCell.Copy()
$q1Sheet.Range("A2:A" + $lastRow).PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 3: populate the "总计" summary sheet: the previous summary rows
# plus a new leading row for 2022-Q1, with the index column renumbered.
# ------------------------------------------------------------------
$totalSheet.Cells.Clear()

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 9, 1.58),
    @("2021-Q4", 28, 5.75),
    @("2021-Q3", 12, 1.25),
    @("2021-Q2", 1, 0.3),
    @("2021-Q1", 8, 4.2),
    @("2020-Q4", 4, 4.45)
)

$totalLastRow = 1 + $totalData.Count

$row = 2
foreach ($item in $totalData) {
    $totalSheet.Cells.Item($row, 1).Value = $row - 2
    $totalSheet.Cells.Item($row, 2).Value = $item[0]
    $totalSheet.Cells.Item($row, 3).Value = $item[1]
    $totalSheet.Cells.Item($row, 4).Value = $item[2]
    $row = $row + 1
}

$refHeaderCell.Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)
$refIndexCell.Copy()
$totalSheet.Range("A2:A" + $totalLastRow).PasteSpecial(-4122)
